# Generate Report for Handback
# Updates the timestamp strings in the handback-status workbook to reflect
# a newly generated report: the "Latest HO Xliff Generate Date" for the
# 6f45d030... file on the Overview sheet, and the corresponding
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" values
# on the zh-cn and de-de sheets for that same file's row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to 6f45d030-ba2f-4e75-aa85-1782b2abefde.md
# (its "Latest HO Xliff Generate Date" shares the same text as de-de's
# "Correspond Handoff Datetime" below, since Excel dedups identical strings)
$wsOverview.Range("G3").Value = "2016-08-18 18:48:13"

# zh-cn sheet: row 3 corresponds to 6f45d030-ba2f-4e75-aa85-1782b2abefde...zh-cn.xlf
$wsZhCn.Range("H3").Value = "2016-08-18 18:48:00"
$wsZhCn.Range("K3").Value = "2016-08-18 18:48:30"

# de-de sheet: row 3 corresponds to 6f45d030-ba2f-4e75-aa85-1782b2abefde...de-de.xlf
$wsDeDe.Range("H3").Value = "2016-08-18 18:48:13"
$wsDeDe.Range("K3").Value = "2016-08-18 18:48:38"
